$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so numeric-looking
# strings (e.g. "1.00", "574.36") are not silently coerced into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "69.415.42"
$ws.Range("E2").Value2 = "  -0.90%  "
$ws.Range("D3").Value2 = "3.516.60"
$ws.Range("E3").Value2 = "  -2.02%  "
$ws.Range("D4").Value2 = "1.00"
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").Value2 = "574.36"
$ws.Range("E5").Value2 = "  -0.69%  "
$ws.Range("D6").Value2 = "185.71"
$ws.Range("E6").Value2 = "  -2.71%  "
$ws.Range("D7").Value2 = "3.502.75"
$ws.Range("E7").Value2 = "  -2.30%  "
$ws.Range("E8").Value2 = "  -2.81%  "
$ws.Range("E9").Value2 = "  +0.09%  "
$ws.Range("D10").Value2 = "0.190"
$ws.Range("E10").Value2 = "  +5.27%  "
$ws.Range("E11").Value2 = "  -2.17%  "
$ws.Range("E12").Value2 = "  -2.77%  "
$ws.Range("D13").Value2 = "0.0000303"
$ws.Range("E13").Value2 = "  -1.37%  "
$ws.Range("D14").Value2 = "9.49"
$ws.Range("E14").Value2 = "  -1.61%  "
$ws.Range("D15").Value2 = "4.075.61"
$ws.Range("E15").Value2 = "  -2.08%  "
$ws.Range("D16").Value2 = "19.38"
$ws.Range("E16").Value2 = "  -2.31%  "
$ws.Range("D17").Value2 = "69.305.93"
$ws.Range("E17").Value2 = "  -0.97%  "
$ws.Range("D18").Value2 = "3.501.48"
$ws.Range("E18").Value2 = "  -2.37%  "
$ws.Range("D19").Value2 = "12.33"
$ws.Range("E19").Value2 = "  -2.56%  "
$ws.Range("D20").Value2 = "0.120"
$ws.Range("E20").Value2 = "  -1.05%  "
$ws.Range("D21").Value2 = "548.03"
$ws.Range("E21").Value2 = "  +15.44%  "
$ws.Range("E22").Value2 = "  -2.72%  "
$ws.Range("D23").Value2 = "18.59"
$ws.Range("E23").Value2 = "  -3.22%  "
$ws.Range("D24").Value2 = "4.96"
$ws.Range("E24").Value2 = "  -1.33%  "
$ws.Range("D25").Value2 = "4.45"
$ws.Range("E25").Value2 = "  +1.81%  "
$ws.Range("D26").Value2 = "94.33"
$ws.Range("E26").Value2 = "  -1.01%  "
$ws.Range("E27").Value2 = "  +3.11%  "
$ws.Range("E28").Value2 = "  -1.66%  "
$ws.Range("D29").Value2 = "9.16"
$ws.Range("E29").Value2 = "  -1.40%  "
$ws.Range("D30").Value2 = "31.85"
$ws.Range("E30").Value2 = "  -1.28%  "
$ws.Range("E31").Value2 = "  -4.86%  "
$ws.Range("D32").Value2 = "12.68"
$ws.Range("E32").Value2 = "  +3.82%  "
$ws.Range("D33").Value2 = "64.79"
$ws.Range("E33").Value2 = "  -2.58%  "
$ws.Range("E34").Value2 = "  -3.98%  "
$ws.Range("D35").Value2 = "545.57"
$ws.Range("E35").Value2 = "  -7.37%  "
$ws.Range("D36").Value2 = "0.405"
$ws.Range("E36").Value2 = "  +2.59%  "
$ws.Range("D37").Value2 = "38.19"
$ws.Range("E37").Value2 = "  -2.09%  "
$ws.Range("B38").Value2 = "Dai"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value2 = "1.00"
$ws.Range("E38").Value2 = "  -0.12%  "
$ws.Range("B39").Value2 = "Fetch.AI"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value2 = "3.06"
$ws.Range("E39").Value2 = "  +7.31%  "
$ws.Range("D40").Value2 = "0.0₃0769"
$ws.Range("E40").Value2 = "  -4.13%  "
$ws.Range("E41").Value2 = "  -2.22%  "
$ws.Range("D42").Value2 = "3.11"
$ws.Range("E42").Value2 = "  -2.47%  "
$ws.Range("E43").Value2 = "  -2.81%  "
$ws.Range("D44").Value2 = "3.304.85"
$ws.Range("E44").Value2 = "  +2.63%  "
$ws.Range("D45").Value2 = "2.99"
$ws.Range("E45").Value2 = "  -2.59%  "
$ws.Range("D46").Value2 = "0.0447"
$ws.Range("E46").Value2 = "  +0.59%  "
$ws.Range("E47").Value2 = "  +2.90%  "
$ws.Range("E48").Value2 = "  -2.27%  "
$ws.Range("D49").Value2 = "8.96"
$ws.Range("E49").Value2 = "  -5.32%  "
$ws.Range("D50").Value2 = "0.997"
$ws.Range("E50").Value2 = "  -0.15%  "
$ws.Range("D51").Value2 = "137.40"
$ws.Range("E51").Value2 = "  +2.64%  "
